$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: a / b / java.rmi.UnmarshalException
$ws.Range("A5").Value = "a"
$ws.Range("B5").Value = "b"
$ws.Range("C5").Value = "java.rmi.UnmarshalException"

# Row 6: -123 / 123525 / formula
$ws.Range("A6").Value = -123
$ws.Range("B6").Value = 123525
$ws.Range("C6").Formula = "=A6+B6"
$ws.Range("C6").NumberFormat = "0.0"

# Row 7: (A7 blank) / reg / java.rmi.UnmarshalException
$ws.Range("B7").Value = "reg"
$ws.Range("C7").Value = "java.rmi.UnmarshalException"

# Column C width (bestFit)
$ws.Columns.Item(3).ColumnWidth = 25.21875

# Update selection to C7
$ws.Range("C7").Select()
